$wb = $excel.ActiveWorkbook

# --- Worksheets ---
$wsPurchase = $wb.Worksheets.Item("Purchase 22-23")

# --- Add the new outstanding row (row 19) to "Purchase 22-23" ---
# Copy the formatting of row 17 (same shape of entry) cell-by-cell so the
# new row 19 cells pick up the identical style ids already used in the sheet.
$wsPurchase.Range("A17").Copy() | Out-Null
$wsPurchase.Range("A19").PasteSpecial(-4122) | Out-Null

$wsPurchase.Range("B17").Copy() | Out-Null
$wsPurchase.Range("B19").PasteSpecial(-4122) | Out-Null

$wsPurchase.Range("C17").Copy() | Out-Null
$wsPurchase.Range("C19").PasteSpecial(-4122) | Out-Null

$wsPurchase.Range("D17").Copy() | Out-Null
$wsPurchase.Range("D19").PasteSpecial(-4122) | Out-Null

$wsPurchase.Range("E17").Copy() | Out-Null
$wsPurchase.Range("E19").PasteSpecial(-4122) | Out-Null

$wsPurchase.Range("F17").Copy() | Out-Null
$wsPurchase.Range("F19").PasteSpecial(-4122) | Out-Null

# Now fill in the actual values for the new row.
$wsPurchase.Range("A19").Value = 8
$wsPurchase.Range("B19").Value = 45252
$wsPurchase.Range("C19").Value = 1200
$wsPurchase.Range("D19").Value = "Sanyo and Sanyo"
$wsPurchase.Range("E19").Value = 11500
$wsPurchase.Range("F19").Formula = "=E19"

# --- Active sheet / selection bookkeeping ---
# "Purchase 22-23" becomes the active/selected tab (was "Sale 22-23"),
# with its cursor left at A20, just below the freshly-added row.
$wsPurchase.Activate()
$wsPurchase.Range("A20").Select() | Out-Null
